# Insert a new weekly record row at row 165, shifting all subsequent rows
# (165-240) down by one, so former row 240's data ends up at row 241.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(165).Insert()

# Populate the newly inserted row 165 with the new record's data.
$ws.Cells.Item(165, 1).Value = 4
$ws.Cells.Item(165, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(165, 3).Value = "Los Lagos"
$ws.Cells.Item(165, 4).Value = 44704
$ws.Cells.Item(165, 5).Value = 10
$ws.Cells.Item(165, 6).Value = "Fruta"
$ws.Cells.Item(165, 7).Value = 100101
$ws.Cells.Item(165, 8).Value = "Berries"
$ws.Cells.Item(165, 9).Value = 100101007
$ws.Cells.Item(165, 10).Value = "Kiwi"
$ws.Cells.Item(165, 11).Value = "Hayward"
$ws.Cells.Item(165, 12).Value = "Segunda"
$ws.Cells.Item(165, 13).Value = 300
$ws.Cells.Item(165, 14).Value = 13000
$ws.Cells.Item(165, 15).Value = 13000
$ws.Cells.Item(165, 16).Value = 13000
$ws.Cells.Item(165, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(165, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(165, 19).Value = 867
$ws.Cells.Item(165, 20).Value = 15
